$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.628
$ws.Range("C14").Value = -13.3476
$ws.Range("C21").Value = -12.19920000000001
$ws.Range("D22").Value = -7.851800000000003
$ws.Range("C23").Value = -12.1814
$ws.Range("D24").Value = -7.451199999999996
$ws.Range("C25").Value = -11.22930000000001
$ws.Range("C26").Value = -12.52250000000001
$ws.Range("D28").Value = -7.984499999999994
$ws.Range("C29").Value = -11.02660000000001
$ws.Range("D36").Value = -7.5532
$ws.Range("D45").Value = -7.289599999999997
$ws.Range("D48").Value = -7.539499999999993
$ws.Range("D49").Value = -7.6328
$ws.Range("D52").Value = -8.095900000000002
$ws.Range("C53").Value = -14.2231
$ws.Range("D53").Value = -8.055799999999998
$ws.Range("D54").Value = -7.899699999999999
$ws.Range("C57").Value = -13.8193
$ws.Range("C59").Value = -12.6794
$ws.Range("C69").Value = -10.73589999999999
$ws.Range("D70").Value = -6.843799999999998
$ws.Range("C79").Value = -11.43020000000001
$ws.Range("C83").Value = -13.9661
$ws.Range("D86").Value = -7.987299999999994
$ws.Range("D87").Value = -8.403599999999992
$ws.Range("D89").Value = -8.085099999999995
$ws.Range("C91").Value = -12.70610000000001
$ws.Range("C93").Value = -10.31339999999999
$ws.Range("D101").Value = -8.128700000000002
$ws.Range("C103").Value = -13.23839999999999
